# feat: add 2022-Q1 data
#
# 1) Insert a new sheet "2022-Q1" (duplicated from "2021-Q4" so it keeps
#    the same column layout / styling) right before the "总计" sheet, and
#    fill it with the per-fund holdings data for the new quarter.
# 2) Insert a new top data row into "总计" with the 2022-Q1 summary
#    figures, pushing the existing rows down one and renumbering their
#    index column.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Force the value to be stored as text (not auto-converted to a
    # number) even when it looks numeric, e.g. "19.02" or "010967".
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" worksheet
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "2022-Q1"

# The source sheet only has 8 data rows; extend formatting down to the 2
# extra rows (10 and 11) we need by copying row 9's formats onto them.
$ws.Range("A9:H9").Copy()
$ws.Range("A10:H11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$fundRows = @(
    @("010967", "博道嘉丰混合A", "19.02", "90.02", "4.89", "0.9301", 4),
    @("010147", "博道嘉兴一年持有期混合", "14.00", "92.27", "5.90", "0.8260", 4),
    @("008467", "博道嘉瑞混合A", "16.14", "77.44", "2.96", "0.4777", 8),
    @("008208", "博道嘉泰回报混合", "16.72", "77.97", "2.72", "0.4548", 8),
    @("008793", "博道嘉元混合A", "9.21", "76.45", "2.56", "0.2358", 9),
    @("010968", "博道嘉丰混合C", "3.58", "90.02", "4.89", "0.1751", 4),
    @("008794", "博道嘉元混合C", "1.58", "76.45", "2.56", "0.0404", 9),
    @("008468", "博道嘉瑞混合C", "1.13", "77.44", "2.96", "0.0334", 8),
    @("011987", "财通资管智选核心回报6个月持有期混合型发起式证券投资基金A", "0.16", "38.14", "1.85", "0.0030", 2),
    @("011988", "财通资管智选核心回报6个月持有期混合型发起式证券投资基金C", "0.01", "38.14", "1.85", "0.0002", 2)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $fund = $fundRows[$i]

    $ws.Cells.Item($r, 1).Value = $i
    Set-TextCell $ws $r 2 $fund[0]
    Set-TextCell $ws $r 3 $fund[1]
    Set-TextCell $ws $r 4 $fund[2]
    Set-TextCell $ws $r 5 $fund[3]
    Set-TextCell $ws $r 6 $fund[4]
    Set-TextCell $ws $r 7 $fund[5]
    $ws.Cells.Item($r, 8).Value = $fund[6]
}

$ws.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Insert the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()

# The insert leaves row 2 blank with mismatched formatting; row 3 (the
# old first data row, now shifted down) still has the correct look, so
# copy its formats back up onto row 2.
$totalWs.Range("A3:D3").Copy()
$totalWs.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalWs.Cells.Item(2, 1).Value = 0
Set-TextCell $totalWs 2 2 "2022-Q1"
$totalWs.Cells.Item(2, 3).Value = 10
$totalWs.Cells.Item(2, 4).Value = 3.18

# The pre-existing rows shifted down by one; renumber their 0-based
# index column (A) to match.
for ($r = 3; $r -le 7; $r++) {
    $totalWs.Cells.Item($r, 1).Value = $r - 2
}
